$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44390
$ws.Range("J2").Value = 55
$ws.Range("K2").Value = 6000
$ws.Range("L2").Value = 6000
$ws.Range("M2").Value = 6000
$ws.Range("P2").Value = 6000

# Row 3
$ws.Range("D3").Value = 44365
$ws.Range("J3").Value = 55
$ws.Range("K3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = 5000
$ws.Range("P3").Value = 5000

# Row 4
$ws.Range("D4").Value = 44957
$ws.Range("J4").Value = 20

# Row 5
$ws.Range("D5").Value = 44497
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 4000
$ws.Range("L5").Value = 4000
$ws.Range("M5").Value = 4000
$ws.Range("P5").Value = 4000

# Row 6
$ws.Range("D6").Value = 44777
$ws.Range("J6").Value = 25
$ws.Range("K6").Value = 5000
$ws.Range("L6").Value = 5000
$ws.Range("M6").Value = 5000
$ws.Range("P6").Value = 5000

# Row 7
$ws.Range("D7").Value = 44504
$ws.Range("J7").Value = 55

# Row 9
$ws.Range("D9").Value = 44966
$ws.Range("J9").Value = 40
$ws.Range("K9").Value = 5000
$ws.Range("L9").Value = 5000
$ws.Range("M9").Value = 5000
$ws.Range("P9").Value = 5000

# Row 10
$ws.Range("D10").Value = 44959
$ws.Range("J10").Value = 40
$ws.Range("K10").Value = 5000
$ws.Range("L10").Value = 5000
$ws.Range("M10").Value = 5000
$ws.Range("P10").Value = 5000

# Row 11
$ws.Range("D11").Value = 44498
$ws.Range("J11").Value = 40

# Row 12
$ws.Range("D12").Value = 44509
$ws.Range("J12").Value = 20
$ws.Range("K12").Value = 4000
$ws.Range("L12").Value = 4000
$ws.Range("M12").Value = 4000
$ws.Range("P12").Value = 4000

# Row 13
$ws.Range("D13").Value = 45163
$ws.Range("J13").Value = 30

# Row 14
$ws.Range("D14").Value = 45169
$ws.Range("J14").Value = 50
$ws.Range("L14").Value = 5000
$ws.Range("M14").Value = 4600
$ws.Range("P14").Value = 4600

# Row 15
$ws.Range("D15").Value = 44956
$ws.Range("K15").Value = 5000
$ws.Range("L15").Value = 5000
$ws.Range("M15").Value = 5000
$ws.Range("P15").Value = 5000

# Row 16
$ws.Range("D16").Value = 44656
$ws.Range("J16").Value = 85
$ws.Range("K16").Value = 5000
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = 5000
$ws.Range("P16").Value = 5000

# Row 17
$ws.Range("D17").Value = 44301
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = 3000
$ws.Range("P17").Value = 3000

# Row 18
$ws.Range("D18").Value = 44679
$ws.Range("J18").Value = 50

# Row 19
$ws.Range("D19").Value = 44313
$ws.Range("J19").Value = 20
$ws.Range("K19").Value = 4000
$ws.Range("L19").Value = 4000
$ws.Range("M19").Value = 4000
$ws.Range("P19").Value = 4000

# Row 20
$ws.Range("D20").Value = 45194
$ws.Range("K20").Value = 6000
$ws.Range("L20").Value = 6000
$ws.Range("M20").Value = 6000
$ws.Range("P20").Value = 6000

# Row 21
$ws.Range("D21").Value = 44680
$ws.Range("J21").Value = 20
$ws.Range("K21").Value = 5000
$ws.Range("M21").Value = 5000
$ws.Range("P21").Value = 5000

# Row 22
$ws.Range("D22").Value = 44259
$ws.Range("J22").Value = 30
$ws.Range("K22").Value = 4000
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = 4000
$ws.Range("P22").Value = 4000

# Row 23
$ws.Range("D23").Value = 44649
$ws.Range("J23").Value = 20
$ws.Range("K23").Value = 5000
$ws.Range("L23").Value = 5000
$ws.Range("M23").Value = 5000
$ws.Range("P23").Value = 5000

# Row 24
$ws.Range("D24").Value = 45159
$ws.Range("J24").Value = 75

# Row 25
$ws.Range("D25").Value = 45162
$ws.Range("J25").Value = 30

# Row 26
$ws.Range("D26").Value = 44781
$ws.Range("J26").Value = 40
$ws.Range("K26").Value = 5000
$ws.Range("L26").Value = 5000
$ws.Range("M26").Value = 5000
$ws.Range("P26").Value = 5000

# Row 27
$ws.Range("D27").Value = 44749
$ws.Range("J27").Value = 65
$ws.Range("K27").Value = 6000
$ws.Range("L27").Value = 6000
$ws.Range("M27").Value = 6000
$ws.Range("P27").Value = 6000

# Row 28
$ws.Range("D28").Value = 44280
$ws.Range("J28").Value = 55
$ws.Range("K28").Value = 4000
$ws.Range("L28").Value = 4000
$ws.Range("M28").Value = 4000
$ws.Range("P28").Value = 4000

# Row 29
$ws.Range("D29").Value = 44316
$ws.Range("J29").Value = 20
$ws.Range("K29").Value = 4000
$ws.Range("L29").Value = 4000
$ws.Range("M29").Value = 4000
$ws.Range("P29").Value = 4000

# Row 30
$ws.Range("D30").Value = 44312
$ws.Range("J30").Value = 50
$ws.Range("K30").Value = 4000
$ws.Range("L30").Value = 4000
$ws.Range("M30").Value = 4000
$ws.Range("P30").Value = 4000

# Row 31
$ws.Range("D31").Value = 44176
$ws.Range("J31").Value = 10

# Row 32
$ws.Range("D32").Value = 44315

# Row 33
$ws.Range("D33").Value = 44508
$ws.Range("J33").Value = 30
$ws.Range("K33").Value = 4000
$ws.Range("L33").Value = 4000
$ws.Range("M33").Value = 4000
$ws.Range("P33").Value = 4000
